$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.430.64'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.940.38'
$ws.Range("E3").Value = '  -2.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.59'
$ws.Range("E5").Value = '  -2.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.73'
$ws.Range("E6").Value = '  +1.48%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.518'
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.936.94'
$ws.Range("E9").Value = '  -2.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.70'
$ws.Range("E10").Value = '  -5.30%  '
$ws.Range("E11").Value = '  -3.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000246'
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.34'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.371.50'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.426.42'
$ws.Range("E17").Value = '  -2.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.94'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.934.94'
$ws.Range("E19").Value = '  -2.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.47'
$ws.Range("E20").Value = '  +11.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '446.12'
$ws.Range("E21").Value = '  -4.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.688'
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.27'
$ws.Range("E23").Value = '  -1.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.20'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.22'
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.11'
$ws.Range("E26").Value = '  -3.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.15'
$ws.Range("E27").Value = '  -5.47%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.06'
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.40'
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("E31").Value = '  -1.72%  '
$ws.Range("E32").Value = '  -4.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.15'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.975'
$ws.Range("E36").Value = '  -2.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.74'
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.74'
$ws.Range("E38").Value = '  +0.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '45.44'
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.98'
$ws.Range("E40").Value = '  -9.33%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.85'
$ws.Range("E41").Value = '  -5.79%  '
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.120'
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.49'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '385.50'
$ws.Range("E45").Value = '  -2.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0352'
$ws.Range("E46").Value = '  -0.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.704.23'
$ws.Range("E47").Value = '  -3.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.22'
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.19'
$ws.Range("E50").Value = '  +4.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.46'
$ws.Range("E51").Value = '  -1.59%  '
